$wb = $excel.ActiveWorkbook
$ws5 = $wb.Worksheets.Item("control_obs")
$ws5.Range("BS20").Formula = "=SUM(BS2:BS18)"
